$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AP3:AP18").Value = 102718.40575454501

$ws.Range("V21").Formula = "=MAX(V2:V20)"
$ws.Range("W21").Formula = "=MIN(W2:W20)"

$ws.Range("I16").Select()
